$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values stay as text (matches source formatting)
$textCells = @("D5", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D18", "D20", "D22", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D43", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.650.51'
$ws.Range("E2").Value = '  -2.32%  '
$ws.Range("D3").Value = '1.789.02'
$ws.Range("E3").Value = '  -2.05%  '
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").Value = '307.84'
$ws.Range("E5").Value = '  -1.75%  '
$ws.Range("E6").Value = '  +0.30%  '
$ws.Range("D7").Value = '0.4558'
$ws.Range("E7").Value = '  +1.70%  '
$ws.Range("D8").Value = '0.3694'
$ws.Range("E8").Value = '  -2.28%  '
$ws.Range("D9").Value = '0.07199'
$ws.Range("E9").Value = '  -4.29%  '
$ws.Range("D10").Value = '0.8535'
$ws.Range("E10").Value = '  -4.47%  '
$ws.Range("D11").Value = '20.37'
$ws.Range("E11").Value = '  -3.05%  '
$ws.Range("D12").Value = '1.795.27'
$ws.Range("E12").Value = '  -1.70%  '
$ws.Range("D13").Value = '5.282'
$ws.Range("E13").Value = '  -2.02%  '
$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").Value = '0.07028'
$ws.Range("E14").Value = '  -1.33%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '6.457'
$ws.Range("E15").Value = '  -4.56%  '
$ws.Range("D16").Value = '90.12'
$ws.Range("E16").Value = '  -4.78%  '
$ws.Range("E17").Value = '  +0.34%  '
$ws.Range("D18").Value = '0.000008595'
$ws.Range("E18").Value = '  -2.30%  '
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("D20").Value = '14.54'
$ws.Range("E20").Value = '  -3.90%  '
$ws.Range("D21").Value = '26.657.83'
$ws.Range("E21").Value = '  -2.29%  '
$ws.Range("D22").Value = '5.267'
$ws.Range("E22").Value = '  +0.43%  '
$ws.Range("E23").Value = '  -3.75%  '
$ws.Range("D24").Value = '2.008.84'
$ws.Range("E24").Value = '  -1.86%  '
$ws.Range("D25").Value = '1.906'
$ws.Range("E25").Value = '  -4.22%  '
$ws.Range("D26").Value = '149.58'
$ws.Range("E26").Value = '  -1.59%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '18.05'
$ws.Range("E27").Value = '  -3.08%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = '2.129'
$ws.Range("E28").Value = '  -13.45%  '
$ws.Range("D29").Value = '5.192'
$ws.Range("E29").Value = '  -3.57%  '
$ws.Range("D30").Value = '113.74'
$ws.Range("E30").Value = '  -3.97%  '
$ws.Range("D31").Value = '0.08814'
$ws.Range("E31").Value = '  -0.42%  '
$ws.Range("D32").Value = '0.7531'
$ws.Range("E32").Value = '  -2.65%  '
$ws.Range("E33").Value = '  -3.14%  '
$ws.Range("D34").Value = '4.429'
$ws.Range("E34").Value = '  -3.36%  '
$ws.Range("D35").Value = '2.882'
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("D36").Value = '1.000'
$ws.Range("E36").Value = '  +0.27%  '
$ws.Range("E37").Value = '  -0.37%  '
$ws.Range("D38").Value = '0.01933'
$ws.Range("D39").Value = '0.05192'
$ws.Range("E39").Value = '  -2.33%  '
$ws.Range("D40").Value = '2.887'
$ws.Range("E40").Value = '  +1.16%  '
$ws.Range("D41").Value = '7.111'
$ws.Range("E41").Value = '  -4.55%  '
$ws.Range("E42").Value = '  +2.92%  '
$ws.Range("D43").Value = '0.5193'
$ws.Range("E43").Value = '  -2.77%  '
$ws.Range("E44").Value = '  -5.41%  '
$ws.Range("D45").Value = '8.445'
$ws.Range("E45").Value = '  -4.05%  '
$ws.Range("D46").Value = '0.4937'
$ws.Range("E46").Value = '  -3.58%  '
$ws.Range("D47").Value = '10.22'
$ws.Range("E47").Value = '  -4.67%  '
$ws.Range("E48").Value = '  +0.21%  '
$ws.Range("D49").Value = '103.76'
$ws.Range("E49").Value = '  -2.36%  '
$ws.Range("D50").Value = '1.639'
$ws.Range("E50").Value = '  -4.10%  '
$ws.Range("D51").Value = '0.06266'
$ws.Range("E51").Value = '  -1.75%  '
